# Update average_county_temperature (column AA) values with NOAA data.
# These rows are grouped by facility/state and each group receives the
# same updated temperature value, replacing the old placeholder of 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AA30:AA33").Value = 19.79629629629628
$ws.Range("AA34:AA37").Value = 16.86342592592595
$ws.Range("AA38:AA41").Value = 5.486111111111112
$ws.Range("AA42:AA57").Value = 14.96875
$ws.Range("AA58:AA61").Value = 17.25771604938272
$ws.Range("AA90:AA105").Value = 13.75752314814816
$ws.Range("AA110:AA121").Value = 14.96875
$ws.Range("AA142:AA157").Value = 13.0158303464755
$ws.Range("AA158:AA161").Value = 16.86342592592595
$ws.Range("AA162:AA169").Value = -3.847222222222223
$ws.Range("AA178:AA181").Value = 12.41429539295394
$ws.Range("AA182:AA185").Value = 19.60879629629628
